# Updated cryptos list on Thu Aug 22 23:35:51 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for rows 2-51 with the
# latest scraped quotes. Values are plain text (as in the source sheet), so
# for Price cells that look like plain decimals (e.g. "584.72") we briefly
# force a Text number format before assigning, then restore the cell's
# default style so no stray formatting is left behind. Price cells that are
# not valid numbers (e.g. "60.327.78", with multiple dots) already stay text
# automatically and don't need the extra step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.327.78'
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").Value = '2.613.51'
$ws.Range("E3").Value = '  -0.55%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.50'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.21%  '

$ws.Range("E10").Value = '  -0.78%  '

$ws.Range("E11").Value = '  +2.28%  '

$ws.Range("E12").Value = '  +1.22%  '

$ws.Range("D13").Value = '3.072.71'
$ws.Range("E13").Value = '  +0.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.17%  '

$ws.Range("D15").Value = '60.307.55'
$ws.Range("E15").Value = '  -0.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000141'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.21%  '

$ws.Range("D17").Value = '2.616.60'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '347.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.07%  '

$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("E23").Value = '  +1.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '

$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.32%  '

$ws.Range("D29").Value = '0.0₃0800'
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.33%  '

$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("E34").Value = '  +6.46%  '

$ws.Range("E35").Value = '  +8.75%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.17%  '

$ws.Range("E37").Value = '  +2.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '319.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.86%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.852'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '135.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0992'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.58%  '

$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.55%  '

$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("E47").Value = '  +0.89%  '

$ws.Range("E48").Value = '  +3.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.83%  '

$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.48%  '
